# Update values in "Name of Algo" result data sheet (terrestrial_mammals / AF / seed3 / KNN)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value  = -21.107
$ws.Range("A10").Value = -20.945
$ws.Range("A12").Value = -21.694
$ws.Range("E13").Value = 12.817
$ws.Range("A18").Value = -21.694
